# Applies the "Updated cryptos list" data refresh to Sheet1 (coin price / 1h volume table).
# Numeric-looking text values (e.g. "497.73") are prefixed with a leading apostrophe so Excel
# keeps storing them as text (quote-prefixed), matching the source data which keeps prices as
# plain strings (some values use "." as a thousands separator, e.g. "56.356.29") rather than numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.356.29"
$ws.Range("E2").Value = "  -4.49%  "
$ws.Range("D3").Value = "2.372.82"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'497.73"
$ws.Range("E5").Value = "  -6.95%  "
$ws.Range("D6").Value = "'128.41"
$ws.Range("E6").Value = "  -4.26%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("E8").Value = "  -3.26%  "
$ws.Range("D9").Value = "2.393.83"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "'0.0954"
$ws.Range("E10").Value = "  -4.81%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  -9.66%  "
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "2.796.68"
$ws.Range("E14").Value = "  -5.23%  "
$ws.Range("D15").Value = "56.235.27"
$ws.Range("E15").Value = "  -4.38%  "
$ws.Range("D16").Value = "'21.35"
$ws.Range("E16").Value = "  -4.49%  "
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("D18").Value = "2.388.90"
$ws.Range("E18").Value = "  -4.61%  "
$ws.Range("D19").Value = "'10.08"
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").Value = "'308.37"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").Value = "'4.01"
$ws.Range("E21").Value = "  -5.45%  "
$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'64.82"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").Value = "2.490.61"
$ws.Range("E26").Value = "  -5.13%  "
$ws.Range("D27").Value = "'0.374"
$ws.Range("E27").Value = "  -8.34%  "
$ws.Range("D28").Value = "'0.150"
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("D29").Value = "'7.20"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").Value = "'172.34"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").Value = "0.0₃0710"
$ws.Range("E31").Value = "  -6.25%  "
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.09"
$ws.Range("E35").Value = "  -7.73%  "
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "'3.77"
$ws.Range("E39").Value = "  -3.96%  "
$ws.Range("D40").Value = "'35.85"
$ws.Range("E40").Value = "  -2.05%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.42"
$ws.Range("E41").Value = "  -6.27%  "
$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D42").Value = "'0.782"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").Value = "'128.80"
$ws.Range("E43").Value = "  -2.29%  "
$ws.Range("E44").Value = "  -4.88%  "
$ws.Range("D45").Value = "'4.77"
$ws.Range("E45").Value = "  -4.89%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'251.69"
$ws.Range("E46").Value = "  -8.52%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.561"
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("E48").Value = "  -4.72%  "
$ws.Range("D49").Value = "'0.0482"
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").Value = "'16.73"
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("E51").Value = "  -5.92%  "
